$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.852.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.743.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5136'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2769'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '38.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06088'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.736.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06988'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6331'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.506'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.870.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006576'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.958.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.089'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.498'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.101'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.40%  '

$ws.Range("E27").Value = '  +2.96%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.817'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08279'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.621'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.397'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04405'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '

$ws.Range("E35").Value = '  -1.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9698'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5981'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.673'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01547'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.907'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3825'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7239'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.874'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05472'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.247'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1101'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '29.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.465'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.48%  '

